$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet/tab to reflect the new "through" date.
$ws.Name = "Through 2022-12-01"

# "November (through 11-30)" -> "November" (A12 keeps referencing the same
# label cell; just update its text).
$ws.Range("A12").Value = "November"

# The old row 13 ("Total") becomes the new row 14; copy its formatting down
# first (while row 13 still holds the original formatting) so the new Total
# row keeps the bold/border/centered style used for the month labels.
$ws.Range("A13:I13").Copy()
$ws.Range("A14:I14").PasteSpecial(-4122)  # xlPasteFormats

# Write the (slightly updated) Total figures into the new row 14.
$ws.Range("A14").Value = "Total"
$ws.Range("B14").Value = 291
$ws.Range("C14").Value = 566
$ws.Range("D14").Value = 825
$ws.Range("E14").Value = 685
$ws.Range("F14").Value = 534
$ws.Range("G14").Value = 1268
$ws.Range("H14").Value = 1649
$ws.Range("I14").Value = 1522

# Repurpose row 13 as the new "December (through 12-01)" data row. 2015 and
# 2019 (columns B and F) have no December figures yet, so those cells are
# cleared rather than set to zero.
$ws.Range("A13").Value = "December (through 12-01)"
$ws.Range("B13").ClearContents()
$ws.Range("C13").Value = 3
$ws.Range("D13").Value = 4
$ws.Range("E13").Value = 3
$ws.Range("F13").ClearContents()
$ws.Range("G13").Value = 4
$ws.Range("H13").Value = 6
$ws.Range("I13").Value = 6
